$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift the base date (B2) forward by 189 days; all other B-column cells
# are formulas relative to B2 (or B4, B6, ...) so they recompute automatically.
$ws.Range("B2").Value = 43374

# Row 20 used to be "Memorial day (no class)" with link_it = FALSE.
# It now becomes the "Building Shiny applications" session (link_it = TRUE).
$ws.Range("D20").Value = "Building Shiny applications"
$ws.Range("C20").Value = $true

# Row 21 used to be "Building Shiny applications"; it becomes part II.
$ws.Range("D21").Value = "Building Shiny applications (part II)"

# Remove the lab rows (22-31) entirely - labs are no longer tracked here.
$ws.Range("A22:D31").EntireRow.Delete()

# Update the active selection to reflect the new last cell used in Excel.
$ws.Range("D22").Select()
